# "Authoring code changes has beeen done"
#
# The "Test Cases" sheet tracked a Result column (D) for each authoring
# test case. The pre-filled "PASS" placeholder results for the
# AuthoringTest / AuthoringAppreciateTest rows (D2:D3) are cleared out,
# leaving the Result column blank (ready for the real run results),
# mirroring the author's selection of D2:D5 before deleting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$resultRange = $ws.Range("D2:D5")
$resultRange.Select()
$resultRange.ClearContents()
